$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete columns U:AD entirely (for rows 1 through current used extent), shrinking the table
#    back to columns A:T before we extend downward with new rows.
$ws.Range("U1:AD29").Delete()

# 2) Row labels in column B (rows 2-29). Row 2 keeps "HKL", rows 3-19 get updated text
#    (this also updates the underlying shared-string table), and rows 20-29 are brand new.
$rowLabels = @(
  "HKL",
  "Spiral5",
  "RotRing OmegaMax-90",
  "Equal Angle",
  "Tilt Rotate",
  "CLR",
  "Rizzie Hex",
  "Thomas Hex",
  "Tilt Rotate_Partial",
  "RotRing OmegaMax-60",
  "Equal Angle_Partial",
  "Rizzie Hex_Partial",
  "ND Single",
  "RD Single",
  "TD Single",
  "Morris Single",
  "Ring Perpendicular to ND",
  "Ring Perpendicular to RD",
  "Ring Perpendicular to TD",
  "OffsetFTD",
  "OffsetATD",
  "OffsetF45",
  "OffsetA45",
  "OffsetFRD",
  "OffsetARD",
  "Gaussian Quadrature",
  "Michael-CCHex",
  "Michael-SNHex"
)

for ($i = 0; $i -lt $rowLabels.Length; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 2).Value = $rowLabels[$i]
}

# 2b) The HKL-index column headers on row 2 (C2:J2) are reordered as part of this edit;
#     K2:T2 (the Pairs headers) are unchanged.
$hklHeaders = @(
  "[1, 1, 0]",
  "[2, 2, 2]",
  "[3, 1, 0]",
  "[3, 2, 1]",
  "[2, 1, 1]",
  "[4, 0, 0]",
  "[2, 2, 0]",
  "[2, 0, 0]"
)
for ($i = 0; $i -lt $hklHeaders.Length; $i++) {
  $c = $i + 3
  $ws.Cells.Item(2, $c).Value = $hklHeaders[$i]
}

# 3) Column A holds the zero-based index (row - 2), with the bold/bordered style already
#    present for rows 2-19; extend that same pattern and style down through row 29 by
#    copying the formatting from the row above (then overwriting the value).
for ($r = 20; $r -le 29; $r++) {
  $srcCell = $ws.Cells.Item($r - 1, 1)
  $dstCell = $ws.Cells.Item($r, 1)
  $srcCell.Copy($dstCell)
  $dstCell.Value = $r - 2
}

# 4) Fill in the data block (columns C:T) with 1s for the newly added rows 20-29,
#    matching the existing rows 3-19.
for ($r = 20; $r -le 29; $r++) {
  for ($c = 3; $c -le 20; $c++) {
    $ws.Cells.Item($r, $c).Value = 1
  }
}
